# ------------------------------------------------------------------
# paises.xlsx refresh: new scrape pass updates case counts for several
# countries; two countries each overtake a neighbour in total cases
# (Japon > Polonia, Belice > Vietnam/Republica del Chad), and the
# "last updated" banner timestamp advances.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Refresh the "Datos actualizados ..." banner in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 01:31"

# 2) Countries that swapped rank (country name moves, data follows the
#    country it now refers to) - set the label for the affected rows
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Polonia"
$ws.Range("A163").Value = "Belice"
$ws.Range("A164").Value = "Vietnam"
$ws.Range("A165").Value = "Republica del Chad"

# 3) Updated statistics: Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes (columns B:H)
$updates = @{
    4 = @(6255554, 39962, 3483850, 2582845, 0, 1123, 188859)
    5 = @(3952790, 41889, 3159096, 671013, 0, 1166, 122681)
    10 = @(624069, 8901, 469557, 134460, 0, 389, 20052)
    13 = @(428239, 10504, 308376, 110944, 0, 259, 8919)
    23 = @(246001, 1209, 221800, 14820, 0, 10, 9381)
    27 = @(129425, 477, 114604, 5689, 0, 6, 9132)
    35 = @(94979, 264, 68736, 24505, 0, 28, 1738)
    43 = @(74893, 819, 62935, 9180, 0, 18, 2778)
    44 = @(71962, 119, 70606, 670, 0, 5, 686)
    47 = @(68392, 527, 57823, 9273, 0, 17, 1296)
    48 = @(67922, 550, 47030, 18834, 0, 19, 2058)
    53 = @(54247, 239, 42010, 11214, 0, 10, 1023)
    55 = @(52440, 468, 49395, 2855, 0, 0, 190)
    74 = @(25117, 499, 18116, 6576, 0, 1, 425)
    89 = @(12381, 284, 11479, 612, 0, 2, 290)
    90 = @(10871, 89, 9348, 1259, 0, 0, 264)
    100 = @(8230, 6, 5828, 2199, 0, 2, 203)
    106 = @(6559, 62, 5241, 1115, 0, 1, 203)
    111 = @(4917, 82, 4058, 757, 0, 2, 102)
    114 = @(4618, 41, 3562, 962, 0, 3, 94)
    147 = @(1797, 38, 685, 1085, 0, 5, 27)
    151 = @(1611, 16, 1419, 148, 0, 0, 44)
    163 = @(1050, 43, 213, 824, 0, 0, 13)
    164 = @(1044, 0, 735, 275, 0, 0, 34)
    165 = @(1017, 4, 884, 56, 0, 0, 77)
    202 = @(29, 1, 23, 4, 0, 0, 2)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $data = New-Object "object[,]" 1,7
    for ($i = 0; $i -lt 7; $i++) {
        $data[0, $i] = $vals[$i]
    }
    $ws.Range("B" + $row + ":H" + $row).Value = $data
}
